# Civ4Archid.xlsx update: refresh terrain/unit docs (AC4-3)
#
# 1. Rename "Ash (Prime Timber)" -> "Yew (Prime Timber)" on the
#    "UnitCombat Weaponry" sheet, and clear the Archer column for the
#    Bronze/Iron weaponry rows (they now only apply to Melee/Mounted).
# 2. Add a new "Terrain" worksheet at the end of the workbook with a
#    terrain-facts reference table (yields, river yields, found/found-coast
#    eligibility, bonus resources and improvement prereqs) plus a couple of
#    cell comments, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. UnitCombat Weaponry sheet updates
# ---------------------------------------------------------------------
$wsUnit = $wb.Worksheets.Item("UnitCombat Weaponry")

# Rename the Prime Timber row label.
$wsUnit.Range("A7").Value = "Yew (Prime Timber)"

# Bronze (Copper) and Iron (Iron) weaponry no longer apply to Archers:
# clear the Archer (column B) values for those two rows, leaving the
# Melee/Mounted (C/D) values untouched.
$wsUnit.Range("B3").ClearContents()
$wsUnit.Range("B4").ClearContents()

# ---------------------------------------------------------------------
# 2. New "Terrain" worksheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTerrain = $wb.Worksheets.Add($null, $lastSheet)
$wsTerrain.Name = "Terrain"

# --- header row: terrain types across the top ---
$wsTerrain.Range("B1").Value = "Grass"
$wsTerrain.Range("C1").Value = "Plains"
$wsTerrain.Range("D1").Value = "Desert"
$wsTerrain.Range("E1").Value = "Tundra"
$wsTerrain.Range("F1").Value = "Snow"
$wsTerrain.Range("G1").Value = "Coast"
$wsTerrain.Range("H1").Value = "Ocean"
$wsTerrain.Range("I1").Value = "Peak"
$wsTerrain.Range("J1").Value = "Hill"

# --- Yield row ---
$wsTerrain.Range("A2").Value = "Yield"
$wsTerrain.Range("B2").Value = "2/0/0"
$wsTerrain.Range("C2").Value = "1/1/0"
$wsTerrain.Range("D2").Value = "-"
$wsTerrain.Range("E2").Value = "-"
$wsTerrain.Range("F2").Value = "-"
$wsTerrain.Range("G2").Value = "1/0/2"
$wsTerrain.Range("H2").Value = "1/0/1"
$wsTerrain.Range("I2").Value = "-"
$wsTerrain.Range("J2").Value = "0/1/0"

# --- RiverYield row ---
$wsTerrain.Range("A3").Value = "RiverYield"
$wsTerrain.Range("B3").Value = "0/0/1"
$wsTerrain.Range("C3").Value = "0/0/1"
$wsTerrain.Range("D3").Value = "0/0/1"
$wsTerrain.Range("E3").Value = "0/0/1"
$wsTerrain.Range("F3").Value = "-"
$wsTerrain.Range("G3").Value = "-"
$wsTerrain.Range("H3").Value = "-"
$wsTerrain.Range("I3").Value = "-"
$wsTerrain.Range("J3").Value = "-"

# --- Found row ---
$wsTerrain.Range("A4").Value = "Found"
$wsTerrain.Range("B4").Value = "x"
$wsTerrain.Range("C4").Value = "x"
$wsTerrain.Range("J4").Value = "x"

# --- Found Coast row ---
$wsTerrain.Range("A5").Value = "Found Coast"
$wsTerrain.Range("E5").Value = "x"

# --- Bonuses row (wrapped, multi-line lists) ---
$wsTerrain.Range("A6").Value = "Bonuses"
$wsTerrain.Range("B6").Value = "Barley`nBison`nCoffee`nCopper`nCorn`nCow`nHorse`nIron`nLead`nObsidian`nPig`nPotato`nRice`nSalt`nSheep`nTea`nUranium`nWheat"
$wsTerrain.Range("C6").Value = "Aluminium`nBarley`nBison`nCopper`nCow`nHorse`nIron`nLead`nObsidian`nPlains`nPotato`nSalt`nSheep`nTobacco`nUranium`nWheat`nWine"
$wsTerrain.Range("D6").Value = "Aluminium`nCopper`nIncense`nIron`nObsidian`nOil`nSalt`nUranium"
$wsTerrain.Range("E6").Value = "Aluminium`nCopper`nDeer`nFur`nHorse`nIron`nLead`nObsidian`nOil`nUranium"
$wsTerrain.Range("F6").Value = "Copper`nFur`nIron`nOil`nUranium"
$wsTerrain.Range("G6").Value = "Clam`nCrab`nFish`nPearls`nShrimp"
$wsTerrain.Range("H6").Value = "Fish`nOil`nShrimp`nWhale"

# --- Improvement/Tech row (wrapped) ---
$wsTerrain.Range("A7").Value = "Improvement/Tech"
$wsTerrain.Range("B7").Value = "Farm/Agriculture`nFort/Mathematics`nHamlet/Sedentary Lifestyle"
$wsTerrain.Range("C7").Value = "Farm/Agriculture`nFort/Mathematics`nHamlet/Sedentary Lifestyle"

# --- comments ---
$wsTerrain.Range("I1").AddComment("Author:`nActually a PlotType in dll") | Out-Null
$wsTerrain.Range("J1").AddComment("Author:`nActually a PlotType in dll") | Out-Null
$wsTerrain.Range("J2").AddComment("Author:`nDefined in the YieldInfos file") | Out-Null

# --- formatting: centre/valign the data, wrap the Bonuses/Improvement rows ---
$dataRange = $wsTerrain.Range("A1:J7")
$dataRange.NumberFormat = "@"
$dataRange.VerticalAlignment = -4108  # xlCenter
$wsTerrain.Range("B1:J5").HorizontalAlignment = -4108  # xlCenter

$wrapRange = $wsTerrain.Range("A6:H7")
$wrapRange.WrapText = $true
$wrapRange.HorizontalAlignment = -4108  # xlCenter

$wsTerrain.Rows(6).RowHeight = 270
$wsTerrain.Rows(7).RowHeight = 60

$wsTerrain.Columns("A").ColumnWidth = 17.67
$wsTerrain.Columns("B").ColumnWidth = 24.17
$wsTerrain.Columns("C").ColumnWidth = 26
$wsTerrain.Range("D1:E1").ColumnWidth = 10

$wsTerrain.PageSetup.PaperSize = 9          # A4
$wsTerrain.PageSetup.Orientation = 1        # xlPortrait

$wsTerrain.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Make the new Terrain sheet the active tab
# ---------------------------------------------------------------------
$wsTerrain.Activate()
